$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14-17: summary labels and stats (bold, size 12, vertically centered)
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Rows("14:14").RowHeight = 15.6

$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$b15 = $ws.Range("B15")
$b15.Font.Bold = $true
$b15.Font.Size = 12
$b15.VerticalAlignment = -4108
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Rows("15:15").RowHeight = 15.6

$ws.Range("B16").Formula = "=MIN(N2:N11)"
$b16 = $ws.Range("B16")
$b16.Font.Bold = $true
$b16.Font.Size = 12
$b16.VerticalAlignment = -4108
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Rows("16:16").RowHeight = 15.6

$ws.Range("B17").Formula = "=MAX(Z2:Z11)"
$b17 = $ws.Range("B17")
$b17.Font.Bold = $true
$b17.Font.Size = 12
$b17.VerticalAlignment = -4108
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Rows("17:17").RowHeight = 15.6

# Row 12: average of column J (|S*|/n) over rows 2-11
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$j12 = $ws.Range("J12")
$j12.Font.Bold = $true
$j12.Font.Size = 11

$ws.Range("A14").Select()
